$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 82.125
$ws.Range("I11").Value = 82.125
$ws.Range("K11").Value = 82.125
$ws.Range("M11").Value = 57.875

# Row 64
$ws.Range("H64").Value = 13161532
$ws.Range("I64").Value = 17860836
$ws.Range("K64").Value = 17860836
$ws.Range("M64").Value = -17860588

# Row 67
$ws.Range("H67").Value = 13161532
$ws.Range("I67").Value = 17860836
$ws.Range("K67").Value = 17860836
$ws.Range("M67").Value = -17859978

# Row 70
$ws.Range("H70").Value = 7536.273

# Row 73
$ws.Range("H73").Value = 7536.273

# Row 88
$ws.Range("H88").Value = 6147.381
$ws.Range("J88").Value = 6568.1577
$ws.Range("L88").Value = 6568.1577
$ws.Range("N88").Value = -7380.1577

# Row 91
$ws.Range("H91").Value = 6147.381
$ws.Range("J91").Value = 6568.1577
$ws.Range("L91").Value = 6568.1577
$ws.Range("N91").Value = -9376.1577

# Row 97
$ws.Range("H97").Value = 5250.25
$ws.Range("J97").Value = 5250.25
$ws.Range("L97").Value = 15750.75
$ws.Range("N97").Value = -16742.75

# Row 101
$ws.Range("H101").Value = 3984.3333
$ws.Range("J101").Value = 3984.3333
$ws.Range("L101").Value = 11952.9999
$ws.Range("N101").Value = -15196.9999

# Row 112
$ws.Range("H112").Value = 119607.35
$ws.Range("J112").Value = 79171.234
$ws.Range("L112").Value = 237513.702
$ws.Range("N112").Value = -239729.702

# Row 113
$ws.Range("H113").Value = 4250
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1903.5385
$ws.Range("I2").Value = 1898.375
$ws.Range("K2").Value = 1898.375
$ws.Range("M2").Value = -1785.375

# Row 32
$ws.Range("H32").Value = 3908.1628
$ws.Range("I32").Value = 1872.7693
$ws.Range("K32").Value = 1872.7693
$ws.Range("M32").Value = -1585.7693

# Row 45
$ws.Range("H45").Value = 6761.1
$ws.Range("J45").Value = 3471
$ws.Range("L45").Value = 3471
$ws.Range("N45").Value = -4225

# Row 97
$ws.Range("H97").Value = 886.0909
$ws.Range("I97").Value = 880.6667
$ws.Range("J97").Value = 910.5
$ws.Range("K97").Value = 880.6667
$ws.Range("L97").Value = 910.5
$ws.Range("M97").Value = -384.6667
$ws.Range("N97").Value = -1902.5

# Row 116
$ws.Range("H116").Value = 1903.5385
$ws.Range("I116").Value = 1898.375
$ws.Range("K116").Value = 1898.375
$ws.Range("M116").Value = 395.625

# Row 122
$ws.Range("H122").Value = 5413.826
$ws.Range("J122").Value = 6784.857
$ws.Range("L122").Value = 20354.571
$ws.Range("N122").Value = -25254.571

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1903.5385
$ws.Range("I3").Value = 1898.375
$ws.Range("K3").Value = 1898.375
$ws.Range("M3").Value = -1784.375

# Row 20
$ws.Range("H20").Value = 2953.8572
$ws.Range("I20").Value = 2853.9285
$ws.Range("J20").Value = 3153.7144
$ws.Range("K20").Value = 2853.9285
$ws.Range("L20").Value = 3153.7144
$ws.Range("M20").Value = -2606.9285
$ws.Range("N20").Value = -3647.7144

# Row 99
$ws.Range("H99").Value = 1503.0769
$ws.Range("I99").Value = 1439.6086
$ws.Range("J99").Value = 1989.6666
$ws.Range("K99").Value = 1439.6086
$ws.Range("L99").Value = 1989.6666
$ws.Range("M99").Value = 58.39139999999998
$ws.Range("N99").Value = -4985.6666

# Row 107
$ws.Range("H107").Value = 66156.375
$ws.Range("I107").Value = 3679.5
$ws.Range("K107").Value = 3679.5
$ws.Range("M107").Value = -1759.5

# Row 134
$ws.Range("H134").Value = 71429630
$ws.Range("I134").Value = 71429630
$ws.Range("K134").Value = 214288890
$ws.Range("M134").Value = -214286355

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1099485.1
$ws.Range("J16").Value = 3296.6667
$ws.Range("L16").Value = 3296.6667
$ws.Range("N16").Value = -3870.6667

# Row 31
$ws.Range("H31").Value = 10735.117
$ws.Range("I31").Value = 6161.8
$ws.Range("J31").Value = 12640.667
$ws.Range("K31").Value = 6161.8
$ws.Range("L31").Value = 12640.667
$ws.Range("M31").Value = -5866.8
$ws.Range("N31").Value = -13230.667

# Row 34
$ws.Range("H34").Value = 10735.117
$ws.Range("I34").Value = 6161.8
$ws.Range("J34").Value = 12640.667
$ws.Range("K34").Value = 6161.8
$ws.Range("L34").Value = 12640.667
$ws.Range("M34").Value = -5959.8
$ws.Range("N34").Value = -13044.667

# Row 97
$ws.Range("H97").Value = 29999.5
$ws.Range("J97").Value = 29999.5
$ws.Range("L97").Value = 29999.5
$ws.Range("N97").Value = -31981.5

# Row 99
$ws.Range("H99").Value = 14865.333
$ws.Range("I99").Value = 15849.875
$ws.Range("K99").Value = 15849.875
$ws.Range("M99").Value = -14351.875

# Row 113
$ws.Range("H113").Value = 1099485.1
$ws.Range("J113").Value = 3296.6667
$ws.Range("L113").Value = 3296.6667
$ws.Range("N113").Value = -7636.6667

# Row 126
$ws.Range("H126").Value = 14865.333
$ws.Range("I126").Value = 15849.875
$ws.Range("K126").Value = 47549.625
$ws.Range("M126").Value = -45079.625

$ws = $wb.Worksheets.Item("CUL")
# Row 51
$ws.Range("H51").Value = 2500
$ws.Range("I51").Value = 2500
$ws.Range("K51").Value = 7500
$ws.Range("M51").Value = -7040

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1470.12
$ws.Range("I97").Value = 1229.9412
$ws.Range("K97").Value = 1229.9412
$ws.Range("M97").Value = -733.9412

# Row 102
$ws.Range("H102").Value = 2868.5908
$ws.Range("I102").Value = 2918.25
$ws.Range("J102").Value = 2372
$ws.Range("K102").Value = 2918.25
$ws.Range("L102").Value = 2372
$ws.Range("M102").Value = -1296.25
$ws.Range("N102").Value = -5616

# Row 107
$ws.Range("H107").Value = 2330.5293
$ws.Range("I107").Value = 1948.4
$ws.Range("K107").Value = 1948.4
$ws.Range("M107").Value = -28.40000000000009

# Row 122
$ws.Range("H122").Value = 117662.09
$ws.Range("I122").Value = 152910.5
$ws.Range("J122").Value = 23666.334
$ws.Range("K122").Value = 458731.5
$ws.Range("L122").Value = 70999.00199999999
$ws.Range("M122").Value = -456281.5
$ws.Range("N122").Value = -75899.00199999999

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3170
$ws.Range("I40").Value = 3170
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3170
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3034
$ws.Range("N40").ClearContents()

# Row 48
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

# Row 61
$ws.Range("H61").Value = 5914.6665
$ws.Range("I61").Value = 5961.4546
$ws.Range("J61").Value = 5400
$ws.Range("K61").Value = 5961.4546
$ws.Range("L61").Value = 5400
$ws.Range("M61").Value = -5759.4546
$ws.Range("N61").Value = -5804

# Row 113
$ws.Range("H113").Value = 5914.6665
$ws.Range("I113").Value = 5961.4546
$ws.Range("J113").Value = 5400
$ws.Range("K113").Value = 5961.4546
$ws.Range("L113").Value = 5400
$ws.Range("M113").Value = -3791.4546
$ws.Range("N113").Value = -9740

$ws = $wb.Worksheets.Item("WVR")
# Row 75
$ws.Range("H75").Value = 123749.75
$ws.Range("J75").Value = 123749.75
$ws.Range("L75").Value = 123749.75
$ws.Range("N75").Value = -125621.75

# Row 78
$ws.Range("H78").Value = 123749.75
$ws.Range("J78").Value = 123749.75
$ws.Range("L78").Value = 371249.25
$ws.Range("N78").Value = -380609.25

# Row 100
$ws.Range("H100").Value = 3677.1428
$ws.Range("I100").Value = 3677.1428
$ws.Range("K100").Value = 7354.2856
$ws.Range("M100").Value = -6813.2856

# Row 113
$ws.Range("H113").Value = 1041.2307
$ws.Range("I113").Value = 878.0833
$ws.Range("K113").Value = 2634.2499
$ws.Range("M113").Value = -464.2498999999998

# Row 132
$ws.Range("H132").Value = 14709492
$ws.Range("I132").Value = 21740946
$ws.Range("J132").Value = 7362
$ws.Range("K132").Value = 65222838
$ws.Range("L132").Value = 22086
$ws.Range("M132").Value = -65220308
$ws.Range("N132").Value = -27146
